$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (shifts rows 6-11 down to 7-12, merges shift automatically).
$ws.Rows.Item(6).Insert()

# Copy the formatting of the row that was just pushed down (now row 7, "Extreme Gradient
# Descent (XGBoost)" row) into the new blank row 6 so the new row matches the table style
# used by non-merged rows (same as rows 3-5).
$ws.Range("A7:C7").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's data.
$ws.Cells.Item(6,1).Value = "Gaussian Naïve Bayes"
$ws.Cells.Item(6,2).Value = "Binary Classification"
$ws.Cells.Item(6,3).Value = "Bank Client Term Deposit"

# Stash the current (pre-rebuild) formatting of the project-name hyperlink column so it
# can be restored after the hyperlinks collection is rebuilt below (adding a hyperlink
# resets the cell to the default "Hyperlink" style).
$ws.Range("C2:C12").Copy()
$ws.Range("Z1:Z11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The hyperlinks collection keeps pointing at the stale (pre-shift) ranges, so rebuild it
# from scratch against the new layout, adding the new row's link as well.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"),  "c. Jupyter Notebooks\Bike Rental Demand.ipynb")      | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"),  "c. Jupyter Notebooks\Temperature Trends.ipynb")      | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"),  "c. Jupyter Notebooks\King County House Sales.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"),  "c. Jupyter Notebooks\Graduate Admissions.ipynb")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"),  "c. Jupyter Notebooks\Bank Client Term Deposit.ipynb")| Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"),  "c. Jupyter Notebooks\Employee Attrition.ipynb")      | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"),  "c. Jupyter Notebooks\Cardiac Risk.ipynb")            | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"),  "c. Jupyter Notebooks\Bank Churn.ipynb")              | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "c. Jupyter Notebooks\Skin Analysis.ipynb")           | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "c. Jupyter Notebooks\Air Passenger.ipynb")           | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "c. Jupyter Notebooks\Champagne Sales.ipynb")         | Out-Null

# Restore the original formatting (font/border/alignment) that Hyperlinks.Add clobbered.
$ws.Range("Z1:Z11").Copy()
$ws.Range("C2:C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1:Z11").Clear()

# Update the active selection to reflect where the user ended up working.
$ws.Range("B8").Select()

$wb.Save()
